# 360MasterData.xlsx - value updates on "1. Scenes Sheet" and "2. Hotspots Sheet"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 2. Hotspots Sheet: refresh the H (X) / I (Y) hotspot-position columns.
# The sheet stores 23 repeating hotspot groups of 7 rows each (rows 2-162);
# every group gets the same new pair of coordinates per offset-in-group.
# ---------------------------------------------------------------------------
$wsHotspots = $wb.Worksheets.Item("2. Hotspots Sheet")
$wsHotspots.Activate()

$hVals = @(5.0199999999999996, 47.13, 2.61, -8.23, 19.239999999999998, -55.99, -17.82)
$iVals = @(28.7, 22.96, -43.25, 173.62, -169.82, -120.79, -41.39)

for ($r = 2; $r -le 162; $r++) {
    $off = ($r - 2) % 7
    $wsHotspots.Cells.Item($r, 8).Value = $hVals[$off]
    $wsHotspots.Cells.Item($r, 9).Value = $iVals[$off]
}

# Leave this sheet's remembered cursor on the last-touched hotspot block.
$wsHotspots.Range("H156:I162").Select()

# ---------------------------------------------------------------------------
# 1. Scenes Sheet: bump the scene's pan/tilt (N2/O2) starting values.
# ---------------------------------------------------------------------------
$wsScenes = $wb.Worksheets.Item("1. Scenes Sheet")
$wsScenes.Activate()

$wsScenes.Range("N2").Value = 14.24
$wsScenes.Range("O2").Value = 75.16

# O2 carries a quote-prefixed number format shared by the whole O column;
# writing .Value resets that format, so restore it from an untouched
# neighbour in the same column before leaving the sheet.
$wsScenes.Range("O3").Copy()
$wsScenes.Range("O2").PasteSpecial(-4122)

# Restore the active selection/cursor shown in the workbook.
$wsScenes.Range("O3").Select()
